$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("10per change")

# --- Fix existing rows 67-69: column E (bsecode) was stored as text, convert to a real number ---
$ws.Range("E67").Value = 509930
$ws.Range("E68").Value = 590024
$ws.Range("E69").Value = 543220

# --- Append three new rows (70-72) carrying the "break out" of the 08:44:49 screener run ---

# Row 70
$ws.Range("A70").Value = "26/06/2024 08:44:49"
$ws.Range("B70").Value = 1
$ws.Range("C70").Value = "SUPREMEIND"
$ws.Range("D70").Value = "Supreme Industries Limited"
$ws.Range("E70").NumberFormat = "@"
$ws.Range("E70").Value = "509930"
$ws.Range("F70").Value = -1.66
$ws.Range("G70").Value = 5790
$ws.Range("H70").Value = 69960

# Row 71
$ws.Range("A71").Value = "26/06/2024 08:44:49"
$ws.Range("B71").Value = 2
$ws.Range("C71").Value = "FACT"
$ws.Range("D71").Value = "Fertilizers And Chemicals Travancore Limited"
$ws.Range("E71").NumberFormat = "@"
$ws.Range("E71").Value = "590024"
$ws.Range("F71").Value = 1.5
$ws.Range("G71").Value = 1014
$ws.Range("H71").Value = 2258833

# Row 72
$ws.Range("A72").Value = "26/06/2024 08:44:49"
$ws.Range("B72").Value = 3
$ws.Range("C72").Value = "MAXHEALTH"
$ws.Range("D72").Value = "Max Healthcare Institute Ltd"
$ws.Range("E72").NumberFormat = "@"
$ws.Range("E72").Value = "543220"
$ws.Range("F72").Value = -1.61
$ws.Range("G72").Value = 879.1
$ws.Range("H72").Value = 1387451
